$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the student row with ID = 2 (Dang, Nhat Huy) - simulating the
# DataGridView "Remove" button which removed the selected student record.
# This is row 3 in the worksheet (row 1 = header, row 2 = ID 1).
$ws.Rows.Item(3).Delete()
